{"js": "// Apply the four wording fixes described in the commit \"Report 2 edits and\n// screen shots.\" Each fix is performed as a targeted search + replace so the\n// rest of the paragraph (and its runs/formatting) stay untouched.\n\nasync function replaceOnce(context, findText, replaceText, options) {\n  const body = context.document.body;\n  const results = body.search(findText, options || { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + findText);\n  }\n\n  results.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// 1. \"tokes\" -> \"tokens\" (typo fix)\nawait replaceOnce(\n  context,\n  \"The placement of tokes must happen in a meticulous order\",\n  \"The placement of tokens must happen in a meticulous order\",\n  { matchCase: true }\n);\n\n// 2. \"two way\" -> \"two ways\"\nawait replaceOnce(\n  context,\n  \"There are two way in which data from the input string of tokens can be derived\",\n  \"There are two ways in which data from the input string of tokens can be derived\",\n  { matchCase: true }\n);\n\n// 3. \"either print out\" -> \"print\"\nawait replaceOnce(\n  context,\n  \"The sole output was to either print out that the input was either accepted or rejected.\",\n  \"The sole output was to print that the input was either accepted or rejected.\",\n  { matchCase: true }\n);\n\n// 4. \"needed changed to\" -> \"needed changes to\"\nawait replaceOnce(\n  context,\n  \"there were needed changed to \",\n  \"there were needed changes to \",\n  { matchCase: true }\n);\n", "ps1": "# Apply the four wording fixes described in the commit \"Report 2 edits and\n# screen shots.\" Each fix is a targeted Find/Replace over the whole document\n# body so the rest of the paragraph (and its runs/formatting) stay untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute([ref]$FindText, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$ReplaceText, 2)\n    if (-not $found) {\n        throw \"Find/Replace failed for: $FindText\"\n    }\n}\n\n# 1. \"tokes\" -> \"tokens\" (typo fix)\nReplace-DocText \"The placement of tokes must happen in a meticulous order\" \"The placement of tokens must happen in a meticulous order\"\n\n# 2. \"two way\" -> \"two ways\"\nReplace-DocText \"There are two way in which data from the input string of tokens can be derived\" \"There are two ways in which data from the input string of tokens can be derived\"\n\n# 3. \"either print out\" -> \"print\"\nReplace-DocText \"The sole output was to either print out that the input was either accepted or rejected.\" \"The sole output was to print that the input was either accepted or rejected.\"\n\n# 4. \"needed changed to\" -> \"needed changes to\"\nReplace-DocText \"there were needed changed to \" \"there were needed changes to \"\n\nWrite-Output \"done\"\n"}
